$excel | Get-Member | Select-String -Pattern "CalculateBeforeSave|SkipCalc"
